{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: insert a new \"Meta description\" paragraph right after the\n//           Heading 1 title paragraph. First run is bold (\"Meta\n//           description\"), second run is the plain-text remainder.\n// Change 2: remove the stray duplicate bold title paragraph near the end\n//           of the document and replace the italic \"Review of 88\n//           Fortunes...\" paragraph's text with the new \"Feature Image\n//           Prompt: ...\" text (keeping the italic formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Change 1: add the \"Meta description\" paragraph after the title ---\nconst titlePara = paragraphs.items[0];\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\n// Reset the new paragraph back to the default (Normal) style instead of\n// inheriting \"Heading 1\" from the title paragraph it was split from.\nmetaPara.styleBuiltIn = Word.BuiltInStyleName.normal;\nawait context.sync();\n\nconst labelRange = metaPara.insertText(\"Meta description\", \"End\");\nlabelRange.font.bold = true;\nawait context.sync();\n\nconst restRange = metaPara.insertText(\n  \": Review of 88 Fortunes, a classic Chinese-themed slot game with a good RTP and interesting customization feature. Play free at top online casinos.\",\n  \"End\"\n);\nrestRange.font.bold = false;\nawait context.sync();\n\n// --- Change 2: drop the duplicated bold title paragraph near the end, and\n//     repurpose the italic paragraph's text for the image prompt ---\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nconst duplicateTitlePara = paragraphs.items[count - 2];\nconst reviewPara = paragraphs.items[count - 1];\n\nduplicateTitlePara.delete();\nawait context.sync();\n\nreviewPara.insertText(\n  'Feature Image Prompt: Create a fun and energetic feature image for the game \"88 Fortunes\". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a golden bowl overflowing with coins, with traditional Chinese lanterns and red and gold elements in the background. The tone of the image should be celebratory and upbeat, reflecting the excitement and potential rewards of the game.',\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Change 1: insert a new \"Meta description\" paragraph right after the\n#           Heading 1 title paragraph. \"Meta description\" is bold, the\n#           remainder of the sentence is plain text.\n# Change 2: remove the stray duplicate bold title paragraph near the end\n#           of the document and replace the italic \"Review of 88\n#           Fortunes...\" paragraph's text with the new \"Feature Image\n#           Prompt: ...\" text (keeping the italic formatting).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: add the \"Meta description\" paragraph after the title ---\n$titlePara = $d.Paragraphs(1)\n[void]$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n# Reset the new paragraph back to the default (Normal) style instead of\n# inheriting \"Heading 1\" from the title paragraph it was split from.\n$metaPara.Range.Style = \"Normal\"\n$metaPara.Range.Text = \"Meta description: Review of 88 Fortunes, a classic Chinese-themed slot game with a good RTP and interesting customization feature. Play free at top online casinos.\"\n\n$boldRange = $d.Content\n[void]$boldRange.Find.Execute(\"Meta description\")\n$boldRange.Font.Bold = 1\n\n# --- Change 2: drop the duplicated bold title paragraph near the end, and\n#     repurpose the italic paragraph's text for the image prompt ---\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs($count - 1)\n$duplicateTitlePara.Range.Delete()\n\n$count2 = $d.Paragraphs.Count\n$reviewPara = $d.Paragraphs($count2)\n$reviewRange = $reviewPara.Range\n# Exclude the trailing paragraph mark so only the visible text is replaced\n# (keeps the paragraph's existing leading empty run + italic run intact).\n$reviewRange.MoveEnd(1, -1) | Out-Null\n$reviewRange.Text = 'Feature Image Prompt: Create a fun and energetic feature image for the game \"88 Fortunes\". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a golden bowl overflowing with coins, with traditional Chinese lanterns and red and gold elements in the background. The tone of the image should be celebratory and upbeat, reflecting the excitement and potential rewards of the game.'\n"}
